$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.861.16"
$ws.Range("E2").Value = "  +0.39%  "
$ws.Range("D3").Value = "1.626.01"
$ws.Range("E3").Value = "  -0.30%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.44"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.31%  "
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.255"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.28%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0631"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.19%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.63"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.80%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0787"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.75%  "
$ws.Range("D12").Value = "1.850.59"
$ws.Range("E12").Value = "  -0.33%  "
$ws.Range("E13").Value = "  -0.42%  "
$ws.Range("D14").Value = "1.626.69"
$ws.Range("E14").Value = "  -0.30%  "
$ws.Range("E15").Value = "  -2.45%  "
$ws.Range("D16").Value = "0.0₃0757"
$ws.Range("E16").Value = "  -0.70%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.55"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.81%  "
$ws.Range("D18").Value = "25.844.72"
$ws.Range("E18").Value = "  +0.29%  "
$ws.Range("E19").Value = "  -0.07%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "192.51"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.16%  "
$ws.Range("E21").Value = "  -1.87%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.95"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.27%  "
$ws.Range("E23").Value = "  -0.73%  "
$ws.Range("E24").Value = "  -1.41%  "
$ws.Range("E25").Value = "  -0.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "142.40"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.33%  "
$ws.Range("E27").Value = "  +0.63%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.85"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.43"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.23%  "
$ws.Range("E30").Value = "  -0.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0497"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.30%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.31"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.72%  "
$ws.Range("E33").Value = "  -0.34%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.58"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.56%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.40"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.92%  "
$ws.Range("E36").Value = "  +0.11%  "
$ws.Range("D37").Value = "1.129.93"
$ws.Range("E37").Value = "  -0.37%  "
$ws.Range("E38").Value = "  +0.67%  "
$ws.Range("E39").Value = "  -1.86%  "
$ws.Range("E40").Value = "  +0.82%  "
$ws.Range("E41").Value = "  -0.17%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "99.43"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.10%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.46"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.30%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.796"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.10%  "
$ws.Range("D45").Value = "1.761.63"
$ws.Range("E45").Value = "  -0.25%  "
$ws.Range("E46").Value = "  -0.90%  "
$ws.Range("E47").Value = "  +1.68%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0528"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +4.12%  "
$ws.Range("E49").Value = "  +2.67%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.415"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.83%  "
$ws.Range("E51").Value = "  +1.33%  "
